$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Insert a new row for period "2509" right after the current last
#    data row (row 21), pushing the signature block (rows 26/27) down
#    to rows 27/28.
# ------------------------------------------------------------------
$ws.Rows("22:22").Insert()

# The freshly inserted row 22 is blank; give it the same formatting
# (borders/number formats) that row 21 currently has -- row 21 is the
# bottom row of the table and carries the closing bottom border.
$ws.Range("B21:J21").Copy()
$ws.Range("B22:J22").PasteSpecial(-4122) | Out-Null

# Row 21 now becomes a regular (interior) row of the table, so give it
# the same formatting as the row above it (row 20).
$ws.Range("B20:J20").Copy()
$ws.Range("B21:J21").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 2. Populate the new data row (same worker, new "Periodo Mora" 2509).
# ------------------------------------------------------------------
$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "1039475729"
$ws.Range("D22").Value = "ANA CATALINA MORENO QUICENO"
$ws.Range("E22").Value = "2509"
$ws.Range("F22").Value = 56940
$ws.Range("G22").Value = 1423500

# ------------------------------------------------------------------
# 3. Update the summary fields to reflect the new period.
# ------------------------------------------------------------------
# VALOR MORA total (was 264770, now includes the new 56940)
$ws.Range("E11").Value = 321710

# Cant. Periodos (was 6, now 7)
$ws.Range("F13").Value = 7
